$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new weekly columns before column B (existing B:E - the weekly
# snapshot columns plus the most-recent-rating column - shift right to E:H)
$ws.Range("B1:D1").EntireColumn.Insert()

# New week headers go in the freshly inserted columns
$ws.Range("C1").Value = "Jun_26"
$ws.Range("D1").Value = "Jun_26"

# Two new analyst/bank rows added to the watch list, with no rating
# history yet (so only the placeholder columns will be populated)
$ws.Range("A28").Value = "Benchmark"
$ws.Range("A29").Value = "Evercore ISI"

$ws.Range("B1").Value = "Jun_27"

# The new week columns have no rating activity yet for the existing
# analysts, so fill them with the same "UN" placeholder used elsewhere
$ws.Range("B2:D27").Value = "UN"
$ws.Range("B28:D28").Value = "UN"
$ws.Range("B29:D29").Value = "UN"

# Keep the report's fixed 8-character column widths across the
# now-wider weekly block of columns (C through H)
$ws.Columns("C").ColumnWidth = 7.166666666666666
$ws.Columns("D").ColumnWidth = 7.166666666666666
$ws.Columns("E").ColumnWidth = 7.166666666666666
$ws.Columns("F").ColumnWidth = 7.166666666666666
$ws.Columns("G").ColumnWidth = 7.166666666666666
$ws.Columns("H").ColumnWidth = 7.166666666666666
